# Update result metrics on several worksheets per the target diff.
$wb = $excel.ActiveWorkbook

# --- arbolu sheet ---
$ws = $wb.Worksheets.Item("arbolu")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 0.6121212121212121
$ws.Range("C2").Value = 0.4119706380575945
$ws.Range("D2").Value = 0.6121212121212121
$ws.Range("E2").Value = 0.4826078971533517
$ws.Range("F2").Value = 0.7741169268441995

# --- bosqueu sheet ---
$ws = $wb.Worksheets.Item("bosqueu")
$ws.Range("A2").Value = 110
$ws.Range("B2").Value = 0.6060606060606061
$ws.Range("C2").Value = 0.5653433616742969
$ws.Range("D2").Value = 0.6060606060606061
$ws.Range("E2").Value = 0.542436518809878
$ws.Range("F2").Value = 0.771043771043771

# --- arbolts sheet ---
$ws = $wb.Worksheets.Item("arbolts")
$ws.Range("B2").Value = 0.8708212942316997
$ws.Range("C2").Value = 1.385562392002157
$ws.Range("D2").Value = 1.177099142809201
$ws.Range("E2").Value = -0.6717157196487737

# --- bosquets sheet ---
$ws = $wb.Worksheets.Item("bosquets")
$ws.Range("A2").Value = 124
$ws.Range("B2").Value = 0.6290313206092344
$ws.Range("C2").Value = 0.6401813750093149
$ws.Range("D2").Value = 0.8001133513504914
$ws.Range("E2").Value = 0.2276051412719259

# --- knnts sheet ---
$ws = $wb.Worksheets.Item("knnts")
$ws.Range("A2").Value = 23
$ws.Range("B2").Value = 0.6465181857472143
$ws.Range("C2").Value = 0.7076411505920056
$ws.Range("D2").Value = 0.8412140931962597
$ws.Range("E2").Value = 0.1462132328768061
